$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Trf"
$ws.Range("C2").Value = "Tfr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4796840000000001
$ws.Range("H2").Value = 1.439052
$ws.Range("I2").Value = 0.00309815167247233
$ws.Range("J2").Value = 0.003108049284200215
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8293556666666667
$ws.Range("N2").Value = 2.488067
$ws.Range("O2").Value = 0.9684897047217467
$ws.Range("P2").Value = 0.9787701502025543
$ws.Range("Q2").Value = 0.3978286436093334
$ws.Range("R2").Value = 3.580457792484001
$ws.Range("S2").Value = 0.003000527998455912
$ws.Range("T2").Value = 0.003042065864733586

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Trf"
$ws.Range("C3").Value = "Tfr2"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4796840000000001
$ws.Range("H3").Value = 1.439052
$ws.Range("I3").Value = 0.00309815167247233
$ws.Range("J3").Value = 0.003108049284200215
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.0269835
$ws.Range("N3").Value = 0.053967
$ws.Range("O3").Value = 0.03151029527825326
$ws.Range("P3").Value = 0.02122984979744567
$ws.Range("Q3").Value = 0.012943553214
$ws.Range("R3").Value = 0.07766131928400001
$ws.Range("S3").Value = 0.0000976236740164173
$ws.Range("T3").Value = 0.00006598341946662908

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Trf"
$ws.Range("C4").Value = "Tfr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 22.35043633333333
$ws.Range("H4").Value = 67.051309
$ws.Range("I4").Value = 0.1443555376176879
$ws.Range("J4").Value = 0.1448167077646516
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8293556666666667
$ws.Range("N4").Value = 2.488067
$ws.Range("O4").Value = 0.9684897047217467
$ws.Range("P4").Value = 0.9787701502025543
$ws.Range("Q4").Value = 18.53646102552256
$ws.Range("R4").Value = 166.828149229703
$ws.Range("S4").Value = 0.1398068520023035
$ws.Range("T4").Value = 0.1417422708106474

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Trf"
$ws.Range("C5").Value = "Tfr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 22.35043633333333
$ws.Range("H5").Value = 67.051309
$ws.Range("I5").Value = 0.1443555376176879
$ws.Range("J5").Value = 0.1448167077646516
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.0269835
$ws.Range("N5").Value = 0.053967
$ws.Range("O5").Value = 0.03151029527825326
$ws.Range("P5").Value = 0.02122984979744567
$ws.Range("Q5").Value = 0.6030929988005
$ws.Range("R5").Value = 3.618557992803
$ws.Range("S5").Value = 0.004548685615384342
$ws.Range("T5").Value = 0.003074436954004137

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Trf"
$ws.Range("C6").Value = "Tfr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.46539633333334
$ws.Range("H6").Value = 190.396189
$ws.Range("I6").Value = 0.409906154456342
$ws.Range("J6").Value = 0.4112156745801394
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8293556666666667
$ws.Range("N6").Value = 2.488067
$ws.Range("O6").Value = 0.9684897047217467
$ws.Range("P6").Value = 0.9787701502025543
$ws.Range("Q6").Value = 52.6353860862959
$ws.Range("R6").Value = 473.7184747766631
$ws.Range("S6").Value = 0.3969898904930493
$ws.Range("T6").Value = 0.4024856275744477

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Trf"
$ws.Range("C7").Value = "Tfr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.46539633333334
$ws.Range("H7").Value = 190.396189
$ws.Range("I7").Value = 0.409906154456342
$ws.Range("J7").Value = 0.4112156745801394
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.0269835
$ws.Range("N7").Value = 0.053967
$ws.Range("O7").Value = 0.03151029527825326
$ws.Range("P7").Value = 0.02122984979744567
$ws.Range("Q7").Value = 1.7125185219605
$ws.Range("R7").Value = 10.275111131763
$ws.Range("S7").Value = 0.01291626396329263
$ws.Range("T7").Value = 0.008730047005691657

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Trf"
$ws.Range("C8").Value = "Tfr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 67.05440800000001
$ws.Range("H8").Value = 201.163224
$ws.Range("I8").Value = 0.4330866284717481
$ws.Range("J8").Value = 0.434470202856191
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8293556666666667
$ws.Range("N8").Value = 2.488067
$ws.Range("O8").Value = 0.9684897047217467
$ws.Range("P8").Value = 0.9787701502025543
$ws.Range("Q8").Value = 55.61195324977868
$ws.Range("R8").Value = 500.507579248008
$ws.Range("S8").Value = 0.4194399409275401
$ws.Range("T8").Value = 0.4252464657080883

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Trf"
$ws.Range("C9").Value = "Tfr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 67.05440800000001
$ws.Range("H9").Value = 201.163224
$ws.Range("I9").Value = 0.4330866284717481
$ws.Range("J9").Value = 0.434470202856191
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.0269835
$ws.Range("N9").Value = 0.053967
$ws.Range("O9").Value = 0.03151029527825326
$ws.Range("P9").Value = 0.02122984979744567
$ws.Range("Q9").Value = 1.809362618268
$ws.Range("R9").Value = 10.856175709608
$ws.Range("S9").Value = 0.01364668754420795
$ws.Range("T9").Value = 0.009223737148102686

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Trf"
$ws.Range("C10").Value = "Tfr2"
$ws.Range("D10").Value = "FAPs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.479164
$ws.Range("H10").Value = 2.958328
$ws.Range("I10").Value = 0.009553527781749779
$ws.Range("J10").Value = 0.006389365514817707
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8293556666666667
$ws.Range("N10").Value = 2.488067
$ws.Range("O10").Value = 0.9684897047217467
$ws.Range("P10").Value = 0.9787701502025543
$ws.Range("Q10").Value = 1.226753045329333
$ws.Range("R10").Value = 7.360518271976
$ws.Range("S10").Value = 0.009252493300397846
$ws.Range("T10").Value = 0.006253720244637147

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Trf"
$ws.Range("C11").Value = "Tfr2"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.479164
$ws.Range("H11").Value = 2.958328
$ws.Range("I11").Value = 0.009553527781749779
$ws.Range("J11").Value = 0.006389365514817707
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.5
$ws.Range("M11").Value = 0.0269835
$ws.Range("N11").Value = 0.053967
$ws.Range("O11").Value = 0.03151029527825326
$ws.Range("P11").Value = 0.02122984979744567
$ws.Range("Q11").Value = 0.039913021794
$ws.Range("R11").Value = 0.159652087176
$ws.Range("S11").Value = 0.0003010344813519314
$ws.Range("T11").Value = 0.000135645270180559
